# Apply GIS lat/long fix: swap header labels for P1/Q1 and update
# P/Q coordinate values with higher-precision lat/long pairs.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mvlv_transformers")

# Swap the header cells: P1 should now reference "longitude" (was Q1's value)
# and Q1 should now reference "latitude" (was P1's value).
$ws.Range("P1").Value = "longitude"
$ws.Range("Q1").Value = "latitude"

# Updated, higher-precision coordinate pairs (longitude in P, latitude in Q)
$coords = @(
    @(5.5643530008095192, 51.621006743843317),
    @(5.5530209761477183, 51.624337024457091),
    @(5.5504170565747302, 51.62112197421213),
    @(5.5494637429464388, 51.623510017791183),
    @(5.5516757109692731, 51.622488359419222),
    @(5.5621971334594589, 51.621278156696327),
    @(5.5556037890367023, 51.624107718459797),
    @(5.5492819656122689, 51.626746665451911),
    @(5.5586866183863224, 51.618183633872547),
    @(5.5566792520004826, 51.627937000342328),
    @(5.5482485553741601, 51.626404317432439),
    @(5.5530852770255033, 51.626830630646523),
    @(5.5548837969885767, 51.619672958373222),
    @(5.5534403477137237, 51.621052806352218),
    @(5.5620234232032573, 51.623321770056393),
    @(5.5572400179952099, 51.621776591367023),
    @(5.5607552792176378, 51.619868274254138)
)

for ($i = 0; $i -lt $coords.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 16).Value = $coords[$i][0]
    $ws.Cells.Item($row, 17).Value = $coords[$i][1]
}

# Update sheet view: scroll position and selection
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 4
$ws.Range("P1:Q18").Select()
